$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in D6 (week 5 start date)
$ws.Range("D6").Value = 43596

# Sprint 4 table (rows 13-19) updates
$ws.Range("E14").Value = 30

$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 15

$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 10

$ws.Range("E17").Value = 40
$ws.Range("F17").Value = 23

$ws.Range("D18").Value = 3

$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 4

# Sprint 5 table (rows 25-30) updates
$ws.Range("D25").Value = 3

$ws.Range("E26").Value = 30

$ws.Range("E28").Value = 15

$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 10

# Update the selected view/cell to match the saved sheet view state
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("H8").Select()
